$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: Module A36 MEEN30020 -> MEEN30140
$ws.Range("A36").Value = "MEEN30140"

# Row 37: Module A37 MEEN20020 -> MEEN30030, Stage D37 2 -> 3, Trimester H37 Aut -> Spr
$ws.Range("A37").Value = "MEEN30030"
$ws.Range("D37").Value = 3
$ws.Range("H37").Value = "Spr"

# Row 51: A51 MEEN40020 -> MEEN40170, College B51 Science -> EngArch, School C51 CompSci -> ElecEng
$ws.Range("A51").Value = "MEEN40170"
$ws.Range("B51").Value = "EngArch"
$ws.Range("C51").Value = "ElecEng"

# Row 52: A52 MEEN40170 -> MEEN40160, School C52 ElecEng -> MechEng
$ws.Range("A52").Value = "MEEN40160"
$ws.Range("C52").Value = "MechEng"

# Row 53: A53 MEEN40160 -> EEEN40730 (new module), School C53 MechEng -> ElecEng
$ws.Range("A53").Value = "EEEN40730"
$ws.Range("C53").Value = "ElecEng"

# Sheet view changes: zoom, topLeftCell, selection
$ws.Application.ActiveWindow.Zoom = 85
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("K24").Select()

# Column width changes: remove custom width on column I (9), set column A width to match old column I width
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
$ws.Columns.Item(1).ColumnWidth = 11.42578125
